{"js": "const body = context.document.body;\n\n// 1) \"...sort the exact same set.\" -> \"...sort the exact same data set.\"\nlet r1 = body.search(\"the exact same set.\", { matchCase: true });\nr1.load(\"items\");\nawait context.sync();\nif (r1.items.length > 0) {\n  r1.items[0].insertText(\"the exact same data set.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"I found merge to 2,119 times faster...\" -> \"I found merge sort was up to 2,119 times faster...\"\nlet r2 = body.search(\"I found merge to 2,119\", { matchCase: true });\nr2.load(\"items\");\nawait context.sync();\nif (r2.items.length > 0) {\n  r2.items[0].insertText(\"I found merge sort was up to 2,119\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) \"...slower than bubble sorts next fastest time...\" -> \"...slower than bubble sort's next fastest time...\"\nlet r3 = body.search(\"bubble sorts next fastest\", { matchCase: true });\nr3.load(\"items\");\nawait context.sync();\nif (r3.items.length > 0) {\n  r3.items[0].insertText(\"bubble sort\\u2019s next fastest\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"...sort the exact same set.\" -> \"...sort the exact same data set.\"\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$rng1.Find.Execute(\"the exact same set.\", $false, $false, $false, $false, $false, $true, 1, $false, \"the exact same data set.\", 1)\n\n# 2) \"I found merge to 2,119 times faster...\" -> \"I found merge sort was up to 2,119 times faster...\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$rng2.Find.Execute(\"I found merge to 2,119\", $false, $false, $false, $false, $false, $true, 1, $false, \"I found merge sort was up to 2,119\", 1)\n\n# 3) \"...slower than bubble sorts next fastest time...\" -> \"...slower than bubble sort's next fastest time...\"\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Replacement.ClearFormatting()\n$rng3.Find.Execute(\"bubble sorts next fastest\", $false, $false, $false, $false, $false, $true, 1, $false, \"bubble sort\u2019s next fastest\", 1)\n"}
